$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-25 Sunday" "2025-05-26 Monday"

Replace-Text "332÷6=55, 2" "679÷7=97, 0"
Replace-Text "575÷5=115, 0" "668÷3=222, 2"
Replace-Text "567÷4=141, 3" "812÷8=101, 4"
Replace-Text "646÷8=80, 6" "215÷4=53, 3"
Replace-Text "172÷3=57, 1" "636÷5=127, 1"

Replace-Text "361÷6=60, 1" "758÷7=108, 2"
Replace-Text "536÷4=134, 0" "684÷4=171, 0"
Replace-Text "523÷5=104, 3" "529÷8=66, 1"
Replace-Text "703÷4=175, 3" "299÷7=42, 5"
Replace-Text "966÷3=322, 0" "933÷8=116, 5"

Replace-Text "839÷6=139, 5" "824÷3=274, 2"
Replace-Text "677÷8=84, 5" "851÷4=212, 3"
Replace-Text "693÷2=346, 1" "975÷7=139, 2"
Replace-Text "823÷9=91, 4" "146÷2=73, 0"
Replace-Text "633÷3=211, 0" "465÷5=93, 0"

Replace-Text "283÷2=141, 1" "499÷4=124, 3"
Replace-Text "406÷4=101, 2" "497÷9=55, 2"
Replace-Text "734÷6=122, 2" "892÷3=297, 1"
Replace-Text "460÷3=153, 1" "538÷4=134, 2"
Replace-Text "817÷4=204, 1" "619÷9=68, 7"

Replace-Text "687÷2=343, 1" "499÷3=166, 1"
Replace-Text "726÷9=80, 6" "353÷5=70, 3"
Replace-Text "972÷4=243, 0" "455÷6=75, 5"
Replace-Text "567÷8=70, 7" "502÷7=71, 5"
Replace-Text "547÷9=60, 7" "522÷6=87, 0"
